# Weekly price-sheet update: a new week's observation is inserted at the
# top of the data block (row 45), pushing all the existing observations
# (old rows 45-135) down by one row (new rows 46-136). The sheet's used
# range grows from A1:R135 to A1:R136 automatically once the row is
# inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 45 (shifts rows 45:135 -> 46:136).
$ws.Rows.Item(45).Insert()

# New observation date: 2021-10-08 (serial 44477), stored at midnight to
# match the rest of column D.
$fecha = Get-Date -Year 2021 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0

$ws.Range("A45").Value2 = 10
$ws.Range("B45").Value2 = "Vega Modelo de Temuco"
$ws.Range("C45").Value2 = "La Araucanía"
$ws.Range("D45").Value2 = $fecha
$ws.Range("E45").Value2 = 9
$ws.Range("F45").Value2 = 100112013
$ws.Range("G45").Value2 = "Alcachofa"
$ws.Range("H45").Value2 = "Madrigal"
$ws.Range("I45").Value2 = "Primera"
$ws.Range("J45").Value2 = 50
$ws.Range("K45").Value2 = 12000
$ws.Range("L45").Value2 = 12000
$ws.Range("M45").Value2 = 12000
$ws.Range("N45").Value2 = "`$/caja 40 unidades"
$ws.Range("O45").Value2 = "Región Metropolitana"
$ws.Range("P45").Value2 = 300
$ws.Range("Q45").Value2 = 40
$ws.Range("R45").Value2 = "Hortaliza"
